# Workbook "BB8.xlsx" - shard roster table edit.
#
# Heywood Jablomee's UTC offset was corrected from 22 to 20, which moves
# his row to the top of that UTC group once the roster table is re-sorted
# by the "UTC" column (ascending) - matching the rest of the already
# UTC-ordered table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the UTC value for the "Heywood Jablomee" row (row 35) before re-sorting.
$ws.Range("D35").Value = 20

# Re-sort the table / data range by the "UTC" column (column D), ascending,
# keeping the header row out of the sorted range.
$lo = $ws.ListObjects.Item("Table1")
$sort = $lo.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("D2:D40"), 0, 1) | Out-Null
$sort.Header = 1
$sort.Apply() | Out-Null

# Leave the selection where the user clicked after sorting.
$ws.Range("C33").Select() | Out-Null
